# Add a new "DOB" row to the IEPD mapping table (eCitation sheet).
# A new row is inserted immediately above the existing "First Name" row
# (originally row 16, now row 17), shifting every row below it down by one.
# The new row carries:
#   A -> "DOB"
#   D -> "/wlq-res-doc:WildlifeLicenseQueryResults/wlq-res-ext:WildlifeLicenseReport/nc:Person/nc:PersonBirthDate/nc:Date"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 16

# Push existing row 16 ("First Name" ...) and everything after it down by one.
$ws.Rows($newRow).Insert()

# Start from a clean slate for the new row (no inherited borders/fills from
# the row that used to be here), then fill in the two cells that matter.
$ws.Range("A" + $newRow + ":D" + $newRow).Clear()
$ws.Range("A" + $newRow).Value = "DOB"
$ws.Range("D" + $newRow).Value = "/wlq-res-doc:WildlifeLicenseQueryResults/wlq-res-ext:WildlifeLicenseReport/nc:Person/nc:PersonBirthDate/nc:Date"

# Match the row height used for this new entry.
$ws.Rows($newRow).RowHeight = 19

# Leave the selection where the editor apparently left it.
$ws.Range("B17").Select()
